# Updated cryptos list — apply new Price (col D) / Volume(1h) (col E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.543.62'
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("D3").Value = '2.673.23'
$ws.Range("E3").Value = '  +3.93%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.44'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.588'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("D9").Value = '2.671.74'
$ws.Range("E9").Value = '  +3.80%  '
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.70'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.22%  '
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.357'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.59'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.82%  '
$ws.Range("E15").Value = '  +3.72%  '
$ws.Range("D16").Value = '63.432.10'
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000145'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("D18").Value = '2.688.00'
$ws.Range("E18").Value = '  +4.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.70%  '
$ws.Range("E20").Value = '  +2.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '339.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.85'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.33%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.70'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.61%  '
$ws.Range("E26").Value = '  +2.01%  '
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '541.84'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +19.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("E32").Value = '  +13.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.99'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.94%  '
$ws.Range("D34").Value = '0.0₃0815'
$ws.Range("E34").Value = '  +2.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '172.92'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.16'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +16.40%  '
$ws.Range("E37").Value = '  +2.58%  '
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.21'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.85'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '174.52'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +11.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.22'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.59%  '
$ws.Range("E44").Value = '  +2.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.33'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0565'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.636'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("E48").Value = '  +3.06%  '
$ws.Range("E49").Value = '  +0.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.74'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.44%  '
